# Apply updated transition-probability matrix values.
# The underlying simulation added more games (see commit message), which
# changed the per-state transition counts and therefore the probabilities
# (count / row-total) stored in this lookup-table sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Af0)
$ws.Range("B2").Value = 0.1943661971830986
$ws.Range("C2").Value = 0.523943661971831
$ws.Range("J2").Value = 0.03943661971830986
$ws.Range("P2").Value = 0.1605633802816901
$ws.Range("S2").Value = 0.08169014084507042

# Row 3 (Af1)
$ws.Range("B3").Value = 0.02061855670103093
$ws.Range("C3").Value = 0.03608247422680412
$ws.Range("J3").Value = 0.04639175257731959
$ws.Range("P3").Value = 0.7268041237113402
$ws.Range("S3").Value = 0.1701030927835052

# Row 4 (Af2)
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2702702702702703

# Row 5 (Af3)
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.09606986899563319
$ws.Range("D6").Value = 0.01310043668122271
$ws.Range("E6").Value = 0.004366812227074236
$ws.Range("F6").Value = 0.08296943231441048
$ws.Range("J6").Value = 0.240174672489083
$ws.Range("O6").Value = 0.03493449781659388
$ws.Range("Q6").Value = 0.1441048034934498
$ws.Range("R6").Value = 0.04803493449781659
$ws.Range("S6").Value = 0.3362445414847162

# Row 7 (Ai1)
$ws.Range("B7").Value = 0.136094674556213
$ws.Range("D7").Value = 0.005917159763313609
$ws.Range("E7").Value = 0.01183431952662722
$ws.Range("F7").Value = 0.02366863905325444
$ws.Range("J7").Value = 0.1893491124260355
$ws.Range("O7").Value = 0.03550295857988166
$ws.Range("Q7").Value = 0.1420118343195266
$ws.Range("R7").Value = 0.05917159763313609
$ws.Range("S7").Value = 0.3964497041420119

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.1027027027027027
$ws.Range("D8").Value = 0.01081081081081081
$ws.Range("F8").Value = 0.06216216216216217
$ws.Range("J8").Value = 0.1189189189189189
$ws.Range("O8").Value = 0.02162162162162162
$ws.Range("Q8").Value = 0.2324324324324324
$ws.Range("R8").Value = 0.06216216216216217
$ws.Range("S8").Value = 0.3891891891891892

# Row 9 (Ai3)
$ws.Range("B9").Value = 0.1017699115044248
$ws.Range("D9").Value = 0.02654867256637168
$ws.Range("F9").Value = 0.05309734513274336
$ws.Range("J9").Value = 0.08849557522123894
$ws.Range("O9").Value = 0.008849557522123894
$ws.Range("Q9").Value = 0.2256637168141593
$ws.Range("R9").Value = 0.05309734513274336
$ws.Range("S9").Value = 0.4424778761061947

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.128419452887538
$ws.Range("D10").Value = 0.02127659574468085
$ws.Range("E10").Value = 0.001519756838905775
$ws.Range("F10").Value = 0.07598784194528875
$ws.Range("J10").Value = 0.1253799392097265
$ws.Range("O10").Value = 0.02051671732522796
$ws.Range("Q10").Value = 0.2401215805471125
$ws.Range("R10").Value = 0.05775075987841945
$ws.Range("S10").Value = 0.3290273556231003

# Row 11 (Bf0)
$ws.Range("G11").Value = 0.1458333333333333
$ws.Range("J11").Value = 0.1006944444444444
$ws.Range("K11").Value = 0.1979166666666667
$ws.Range("L11").Value = 0.5486111111111112
$ws.Range("S11").Value = 0.006944444444444444

# Row 12 (Bf1)
$ws.Range("G12").Value = 0.7134146341463414
$ws.Range("J12").Value = 0.2195121951219512
$ws.Range("K12").Value = 0.006097560975609756
$ws.Range("L12").Value = 0.02439024390243903
$ws.Range("S12").Value = 0.03658536585365853

# Row 13 (Bf2)
$ws.Range("G13").Value = 0.5185185185185185
$ws.Range("J13").Value = 0.3703703703703703
$ws.Range("S13").Value = 0.1111111111111111

# Row 15 (Bi0)
$ws.Range("F15").Value = 0.01769911504424779
$ws.Range("H15").Value = 0.1415929203539823
$ws.Range("I15").Value = 0.09292035398230089
$ws.Range("J15").Value = 0.3539823008849557
$ws.Range("K15").Value = 0.04424778761061947
$ws.Range("M15").Value = 0.008849557522123894
$ws.Range("O15").Value = 0.0752212389380531
$ws.Range("S15").Value = 0.2654867256637168

# Row 16 (Bi1)
$ws.Range("F16").Value = 0.01818181818181818
$ws.Range("H16").Value = 0.1363636363636364
$ws.Range("I16").Value = 0.1045454545454545
$ws.Range("J16").Value = 0.4454545454545454
$ws.Range("K16").Value = 0.08181818181818182
$ws.Range("M16").Value = 0.01818181818181818
$ws.Range("N16").Value = 0.004545454545454545
$ws.Range("O16").Value = 0.05454545454545454
$ws.Range("S16").Value = 0.1363636363636364

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.02371541501976284
$ws.Range("H17").Value = 0.1640316205533597
$ws.Range("I17").Value = 0.1225296442687747
$ws.Range("J17").Value = 0.4308300395256917
$ws.Range("K17").Value = 0.07905138339920949
$ws.Range("M17").Value = 0.005928853754940711
$ws.Range("O17").Value = 0.05731225296442688
$ws.Range("S17").Value = 0.116600790513834

# Row 18 (Bi3)
$ws.Range("F18").Value = 0.03816793893129771
$ws.Range("H18").Value = 0.1755725190839695
$ws.Range("I18").Value = 0.0916030534351145
$ws.Range("J18").Value = 0.4427480916030535
$ws.Range("K18").Value = 0.06106870229007633
$ws.Range("M18").Value = 0.02290076335877863
$ws.Range("O18").Value = 0.06106870229007633
$ws.Range("S18").Value = 0.1068702290076336

# Row 19 (Br0)
$ws.Range("F19").Value = 0.01998334721065778
$ws.Range("H19").Value = 0.1698584512905912
$ws.Range("I19").Value = 0.09159034138218151
$ws.Range("J19").Value = 0.3830141548709409
$ws.Range("K19").Value = 0.1223980016652789
$ws.Range("M19").Value = 0.01582014987510408
$ws.Range("N19").Value = 0.001665278934221482
$ws.Range("O19").Value = 0.06994171523730225
$ws.Range("S19").Value = 0.1257285595337219

Write-Output "Updated 113 cells in the transition matrix"
